$wb = $excel.ActiveWorkbook

# --- "readme" sheet: reorder JobNo / sheet_name / Author -> Author / JobNo / sheet_name ---
$ws = $wb.Worksheets.Item("readme")

# Update header row (this also renames the backing table/ListObject columns)
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "JobNo"
$ws.Range("D1").Value = "sheet_name"

# Update data rows 2-12 for columns B (Author), C (JobNo), D (sheet_name)
$ws.Range("B2").Value = "jovyan"
$ws.Range("C2").Value = "Project Information"
$ws.Range("D2").Value = "jovyan"

$ws.Range("B3").Value = "jovyan"
$ws.Range("C3").Value = "/c/e"
$ws.Range("D3").Value = "Criterion Definitions"

$ws.Range("B4").Value = "jovyan"
$ws.Range("C4").Value = "/c/e"
$ws.Range("D4").Value = "Results, Air Speed 0.1"

$ws.Range("B5").Value = "jovyan"
$ws.Range("C5").Value = "/c/e"
$ws.Range("D5").Value = "Results, Air Speed 0.15"

$ws.Range("B6").Value = "jovyan"
$ws.Range("C6").Value = "/c/e"
$ws.Range("D6").Value = "Results, Air Speed 0.2"

$ws.Range("B7").Value = "jovyan"
$ws.Range("C7").Value = "/c/e"
$ws.Range("D7").Value = "Results, Air Speed 0.3"

$ws.Range("B8").Value = "jovyan"
$ws.Range("C8").Value = "/c/e"
$ws.Range("D8").Value = "Results, Air Speed 0.4"

$ws.Range("B9").Value = "jovyan"
$ws.Range("C9").Value = "/c/e"
$ws.Range("D9").Value = "Results, Air Speed 0.5"

$ws.Range("B10").Value = "jovyan"
$ws.Range("C10").Value = "/c/e"
$ws.Range("D10").Value = "Results, Air Speed 0.6"

$ws.Range("B11").Value = "jovyan"
$ws.Range("C11").Value = "Results, Air Speed 0.7"
$ws.Range("D11").Value = "Results, Air Speed 0.7"

$ws.Range("B12").Value = "jovyan"
$ws.Range("C12").Value = "/c/e"
$ws.Range("D12").Value = "Results, Air Speed 0.8"

# --- "Project Information" sheet: bump the analysis timestamp ---
$wsInfo = $wb.Worksheets.Item("Project Information")
$wsInfo.Range("B12").Value = "2022-06-15 15:57:07.449572"
